# Refresh the live cryptocurrency Price (column D) and Volume(1h)
# (column E) figures with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain-text figures (e.g. "45.58", "  +2.68%  ").
# Some of the new Price values parse as plain numbers (e.g. "45.58"),
# so mark those cells as Text first to stop Excel from silently
# converting them to numeric values (which would also introduce
# floating point rounding noise and drop significant trailing zeros).
$textFormatCells = @("D5","D7","D12","D14","D18","D20","D21","D22","D23","D24","D25","D27","D28","D31","D33","D35","D36","D38","D39","D40","D42","D43","D44","D47","D50","D51")
foreach ($ref in $textFormatCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "96.874.27"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.681.36"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "239.63"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  +10.51%  "
$ws.Range("D7").Value = "658.02"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "3.678.58"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "45.58"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D14").Value = "6.81"
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("D15").Value = "4.366.46"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "96.657.42"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "8.92"
$ws.Range("E18").Value = "  +13.03%  "
$ws.Range("D19").Value = "3.673.85"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "18.78"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").Value = "12.75"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "0.530"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "531.40"
$ws.Range("E23").Value = "  +3.21%  "
$ws.Range("D24").Value = "3.51"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "7.15"
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "102.64"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "13.52"
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").Value = "3.04"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "1.91"
$ws.Range("E33").Value = "  +15.93%  "
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").Value = "663.09"
$ws.Range("E35").Value = "  +6.23%  "
$ws.Range("D36").Value = "32.67"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "0.598"
$ws.Range("E38").Value = "  +5.33%  "
$ws.Range("D39").Value = "8.89"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "0.161"
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "6.60"
$ws.Range("E42").Value = "  +9.67%  "
$ws.Range("D43").Value = "0.964"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").Value = "38.76"
$ws.Range("E44").Value = "  +17.26%  "
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("D47").Value = "0.430"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +7.55%  "
$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "8.68"
$ws.Range("E51").Value = "  +1.97%  "
